# InputData/plcy-schd/IT/Initial Time.xlsx
# Adds the "Initial" date value next to the "IT Initial Time" header on the
# About sheet: cell C1 gets a date (2021-04-21, serial 44307) formatted
# with the built-in short-date number format (numFmtId 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Set the number format first so the engine doesn't mint a throwaway
# custom date/time format before we pin it to the built-in "mm-dd-yy"
# (numFmtId 14) short-date format.
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = (Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0)
